$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "45.358.70"
$ws.Range("E2").Value = "  +5.44%  "

# Row 3
$ws.Range("D3").Value = "2.366.92"
$ws.Range("E3").Value = "  +2.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "110.41"
$ws.Range("E5").Value = "  +4.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.78"
$ws.Range("E6").Value = "  -0.75%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +0.54%  "

# Row 8
$ws.Range("E8").Value = "  -0.26%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  +2.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.52"
$ws.Range("E10").Value = "  +3.54%  "

# Row 11
$ws.Range("E11").Value = "  +0.61%  "

# Row 12
$ws.Range("E12").Value = "  +1.64%  "

# Row 13
$ws.Range("E13").Value = "  +1.40%  "

# Row 14
$ws.Range("E14").Value = "  -0.48%  "

# Row 15
$ws.Range("D15").Value = "2.728.11"
$ws.Range("E15").Value = "  +2.53%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.45"
$ws.Range("E16").Value = "  +0.97%  "

# Row 17
$ws.Range("D17").Value = "2.366.43"
$ws.Range("E17").Value = "  +2.43%  "

# Row 18
$ws.Range("D18").Value = "45.339.33"
$ws.Range("E18").Value = "  +5.76%  "

# Row 19
$ws.Range("E19").Value = "  -2.14%  "

# Row 20
$ws.Range("E20").Value = "  +1.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.06"
$ws.Range("E21").Value = "  -2.29%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.62"
$ws.Range("E22").Value = "  +0.08%  "

# Row 23
$ws.Range("E23").Value = "  -1.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.94"
$ws.Range("E24").Value = "  -1.88%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("E25").Value = "  +2.23%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.10"
$ws.Range("E27").Value = "  +1.27%  "

# Row 28
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.41"
$ws.Range("E28").Value = "  -4.60%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  -0.55%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.48"
$ws.Range("E30").Value = "  +0.80%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0967"
$ws.Range("E31").Value = "  +11.27%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.93"
$ws.Range("E32").Value = "  -1.43%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "170.21"
$ws.Range("E33").Value = "  +2.78%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.93"
$ws.Range("E34").Value = "  +5.51%  "

# Row 35
$ws.Range("E35").Value = "  +0.35%  "

# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.83"
$ws.Range("E36").Value = "  +3.97%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.116"
$ws.Range("E37").Value = "  +4.44%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.99"
$ws.Range("E38").Value = "  +6.21%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0357"
$ws.Range("E39").Value = "  -0.24%  "

# Row 40
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.91"
$ws.Range("E40").Value = "  +8.17%  "

# Row 41
$ws.Range("E41").Value = "  +9.26%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.07"
$ws.Range("E42").Value = "  -2.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.234"
$ws.Range("E43").Value = "  +1.71%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.06"
$ws.Range("E44").Value = "  -1.08%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.10"
$ws.Range("E45").Value = "  +6.80%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  -0.46%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.50"
$ws.Range("E47").Value = "  +8.25%  "

# Row 48
$ws.Range("E48").Value = "  +6.19%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.99"
$ws.Range("E49").Value = "  +1.20%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.55"
$ws.Range("E50").Value = "  +6.60%  "

# Row 51
$ws.Range("E51").Value = "  +6.99%  "
